$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text could be misinterpreted as a number by Excel
# (e.g. "1.040", "0.3800", "0.000009151") need to be forced to Text
# format before assignment, then restored to the default "Normal"
# style so no stray formatting is introduced.
function Set-TextValue($ws, $addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

$ws.Range('D2').Value = '27.768.44'
$ws.Range('E2').Value = '  +3.01%  '
$ws.Range('D3').Value = '1.864.06'
$ws.Range('E3').Value = '  +2.74%  '
Set-TextValue $ws 'D4' '1.040'
$ws.Range('E4').Value = '  +3.23%  '
Set-TextValue $ws 'D5' '324.32'
$ws.Range('E5').Value = '  +3.87%  '
Set-TextValue $ws 'D6' '1.036'
$ws.Range('E6').Value = '  +2.99%  '
Set-TextValue $ws 'D7' '0.4427'
$ws.Range('E7').Value = '  +3.05%  '
Set-TextValue $ws 'D8' '0.3800'
$ws.Range('E8').Value = '  +3.51%  '
Set-TextValue $ws 'D9' '0.07463'
$ws.Range('E9').Value = '  +3.20%  '
Set-TextValue $ws 'D10' '0.8848'
$ws.Range('E10').Value = '  +2.54%  '
Set-TextValue $ws 'D11' '21.69'
$ws.Range('E11').Value = '  +2.27%  '
$ws.Range('D12').Value = '1.872.49'
$ws.Range('E12').Value = '  -16.13%  '
Set-TextValue $ws 'D13' '5.562'
$ws.Range('E13').Value = '  +3.04%  '
Set-TextValue $ws 'D14' '6.766'
$ws.Range('E14').Value = '  +2.53%  '
Set-TextValue $ws 'D15' '0.07229'
$ws.Range('E15').Value = '  +4.07%  '
Set-TextValue $ws 'D16' '83.78'
$ws.Range('E16').Value = '  +3.26%  '
Set-TextValue $ws 'D17' '1.040'
$ws.Range('E17').Value = '  +2.78%  '
Set-TextValue $ws 'D18' '0.000009151'
$ws.Range('E18').Value = '  +3.00%  '
Set-TextValue $ws 'D19' '1.036'
$ws.Range('E19').Value = '  +3.02%  '
$ws.Range('E20').Value = '  +2.57%  '
$ws.Range('D21').Value = '27.761.40'
$ws.Range('E21').Value = '  +2.80%  '
Set-TextValue $ws 'D22' '5.319'
$ws.Range('E22').Value = '  +2.97%  '
$ws.Range('E23').Value = '  +3.26%  '
Set-TextValue $ws 'D24' '2.001'
$ws.Range('E24').Value = '  +7.28%  '
Set-TextValue $ws 'D25' '159.01'
$ws.Range('E25').Value = '  +3.33%  '
Set-TextValue $ws 'D26' '18.87'
$ws.Range('E26').Value = '  +3.02%  '
Set-TextValue $ws 'D27' '5.325'
$ws.Range('E27').Value = '  +1.95%  '
$ws.Range('E28').Value = '  +4.22%  '
Set-TextValue $ws 'D29' '118.08'
$ws.Range('E29').Value = '  +3.07%  '
Set-TextValue $ws 'D30' '0.09068'
$ws.Range('E30').Value = '  +1.37%  '
Set-TextValue $ws 'D31' '0.7782'
$ws.Range('E31').Value = '  +3.97%  '
Set-TextValue $ws 'D32' '3.096'
$ws.Range('E32').Value = '  +10.40%  '
$ws.Range('E33').Value = '  +2.02%  '
$ws.Range('E34').Value = '  +3.50%  '
Set-TextValue $ws 'D35' '1.038'
$ws.Range('E35').Value = '  +3.17%  '
Set-TextValue $ws 'D36' '1.153'
$ws.Range('E36').Value = '  +2.61%  '
Set-TextValue $ws 'D37' '0.01993'
$ws.Range('E37').Value = '  +3.79%  '
Set-TextValue $ws 'D38' '0.05347'
$ws.Range('E38').Value = '  +2.65%  '
Set-TextValue $ws 'D39' '2.858'
$ws.Range('E39').Value = '  +3.95%  '
Set-TextValue $ws 'D40' '0.5200'
$ws.Range('E40').Value = '  +1.90%  '
Set-TextValue $ws 'D41' '0.1692'
$ws.Range('E41').Value = '  +2.48%  '
Set-TextValue $ws 'D42' '6.895'
Set-TextValue $ws 'D43' '8.671'
$ws.Range('E43').Value = '  +4.12%  '
Set-TextValue $ws 'D44' '109.82'
$ws.Range('E44').Value = '  +2.85%  '
Set-TextValue $ws 'D45' '10.74'
$ws.Range('E45').Value = '  +3.37%  '
Set-TextValue $ws 'D46' '1.721'
$ws.Range('E46').Value = '  +4.77%  '
Set-TextValue $ws 'D47' '0.4708'
$ws.Range('E47').Value = '  +2.85%  '
Set-TextValue $ws 'D48' '0.06476'
Set-TextValue $ws 'D49' '1.925'
$ws.Range('E49').Value = '  +4.11%  '
Set-TextValue $ws 'D50' '39.88'
$ws.Range('E50').Value = '  +1.95%  '
Set-TextValue $ws 'D51' '64.57'
$ws.Range('E51').Value = '  +2.60%  '
